# Update automàtic: dades i banners [2026-02-22 21:19]
# Refresh DATA_EXTRACCIO timestamps and re-scraped meteo.cat readings
# (humidity/pressure/temperature/radiation) for resum_diari_meteocat.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-22 21:18:19"
$ws.Range("O2").Value = "6.1 °C"
$ws.Range("E3").Value = "2026-02-22 21:18:21"
$ws.Range("N3").Value = "1.4 °C 20:54 TU"
$ws.Range("E4").Value = "2026-02-22 21:18:23"
$ws.Range("E5").Value = "2026-02-22 21:18:26"
$ws.Range("E6").Value = "2026-02-22 21:18:28"
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "62%"
$ws.Range("H6").NumberFormat = "general"
$ws.Range("E7").Value = "2026-02-22 21:18:31"
$ws.Range("E8").Value = "2026-02-22 21:18:33"
$ws.Range("E9").Value = "2026-02-22 21:18:35"
$ws.Range("E10").Value = "2026-02-22 21:18:36"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "81%"
$ws.Range("H10").NumberFormat = "general"
$ws.Range("O10").Value = "10.0 °C"
$ws.Range("E11").Value = "2026-02-22 21:18:37"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "65%"
$ws.Range("H11").NumberFormat = "general"
$ws.Range("O11").Value = "8.7 °C"
$ws.Range("E12").Value = "2026-02-22 21:18:39"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "87%"
$ws.Range("H12").NumberFormat = "general"
$ws.Range("E13").Value = "2026-02-22 21:18:40"
$ws.Range("E14").Value = "2026-02-22 21:18:41"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "74%"
$ws.Range("H14").NumberFormat = "general"
$ws.Range("O14").Value = "11.9 °C"
$ws.Range("E15").Value = "2026-02-22 21:18:42"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "73%"
$ws.Range("H15").NumberFormat = "general"
$ws.Range("E16").Value = "2026-02-22 21:18:43"
$ws.Range("O16").Value = "5.3 °C"
$ws.Range("E17").Value = "2026-02-22 21:18:44"
$ws.Range("O17").Value = "10.0 °C"
$ws.Range("E18").Value = "2026-02-22 21:18:45"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "76%"
$ws.Range("H18").NumberFormat = "general"
$ws.Range("O18").Value = "10.1 °C"
$ws.Range("E19").Value = "2026-02-22 21:18:46"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "47%"
$ws.Range("H19").NumberFormat = "general"
$ws.Range("O19").Value = "12.1 °C"
$ws.Range("E20").Value = "2026-02-22 21:18:47"
$ws.Range("E21").Value = "2026-02-22 21:18:50"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "59%"
$ws.Range("H21").NumberFormat = "general"
$ws.Range("J21").Value = "1028.9 hPa"
$ws.Range("E22").Value = "2026-02-22 21:18:52"
$ws.Range("O22").Value = "4.5 °C"
$ws.Range("E23").Value = "2026-02-22 21:18:54"
$ws.Range("K23").Value = "15.7 MJ/m2"
$ws.Range("O23").Value = "5.8 °C"
$ws.Range("E24").Value = "2026-02-22 21:18:57"
$ws.Range("K24").Value = "15.5 MJ/m2"
$ws.Range("O24").Value = "7.7 °C"
$ws.Range("E25").Value = "2026-02-22 21:18:59"
$ws.Range("E26").Value = "2026-02-22 21:19:01"
$ws.Range("H26").NumberFormat = "@"
$ws.Range("H26").Value = "37%"
$ws.Range("H26").NumberFormat = "general"
$ws.Range("O26").Value = "11.2 °C"
$ws.Range("E27").Value = "2026-02-22 21:19:04"
$ws.Range("K27").Value = "16.2 MJ/m2"
$ws.Range("E28").Value = "2026-02-22 21:19:06"
$ws.Range("E29").Value = "2026-02-22 21:19:09"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "83%"
$ws.Range("H29").NumberFormat = "general"
$ws.Range("O29").Value = "9.8 °C"
$ws.Range("E30").Value = "2026-02-22 21:19:11"
$ws.Range("O30").Value = "12.2 °C"
$ws.Range("E31").Value = "2026-02-22 21:19:14"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "61%"
$ws.Range("H31").NumberFormat = "general"
$ws.Range("E32").Value = "2026-02-22 21:19:16"
$ws.Range("O32").Value = "6.0 °C"
$ws.Range("E33").Value = "2026-02-22 21:19:19"
$ws.Range("O33").Value = "8.3 °C"
$ws.Range("E34").Value = "2026-02-22 21:19:21"
$ws.Range("E35").Value = "2026-02-22 21:19:24"
$ws.Range("E36").Value = "2026-02-22 21:19:26"
$ws.Range("J36").Value = "1027.3 hPa"
$ws.Range("E37").Value = "2026-02-22 21:19:29"
$ws.Range("O37").Value = "8.0 °C"
$ws.Range("E38").Value = "2026-02-22 21:19:31"
$ws.Range("O38").Value = "11.5 °C"
$ws.Range("E39").Value = "2026-02-22 21:19:34"
$ws.Range("E40").Value = "2026-02-22 21:19:36"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "56%"
$ws.Range("H40").NumberFormat = "general"
$ws.Range("O40").Value = "10.0 °C"
$ws.Range("E41").Value = "2026-02-22 21:19:39"
$ws.Range("E42").Value = "2026-02-22 21:19:41"
$ws.Range("E43").Value = "2026-02-22 21:19:43"
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "71%"
$ws.Range("H43").NumberFormat = "general"
$ws.Range("E44").Value = "2026-02-22 21:19:45"
$ws.Range("O44").Value = "2.6 °C"
$ws.Range("E45").Value = "2026-02-22 21:19:48"
$ws.Range("H45").NumberFormat = "@"
$ws.Range("H45").Value = "55%"
$ws.Range("H45").NumberFormat = "general"
$ws.Range("O45").Value = "8.7 °C"
$ws.Range("E46").Value = "2026-02-22 21:19:50"
$ws.Range("J46").Value = "1029.5 hPa"
